# Add mri_convertOptions as a new standard parameter row in the
# functionalParameterProperties sheet.
#
# The new parameter is inserted as a new row 29 (pushing the existing
# row 29 "functional_preprocessing.preprocessingScript" and everything
# below it down by one row), directly after the existing
# "functional_preprocessing.sliceTimerOptions" row (row 28).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 29; existing rows 29+ shift down to 30+.
$ws.Rows("29:29").Insert() | Out-Null

# Populate the new row with the mri_convertOptions parameter definition.
$ws.Range("A29").Value2 = "functional_preprocessing.mri_convertOptions"
$ws.Range("D29").Value2 = "functional_preprocessing"
$ws.Range("E29").Value2 = "char"
$ws.Range("G29").Value2 = "standard"
$ws.Range("H29").Value2 = 'Adjust variables in the fmriProcessedFile header (using mri_convert). Options are provided as structure (e.g. mri_convertOptions:{"tr": TR in msec, "te": TE in msec}). If emtpy, header is not changed.'

# Match the saved selection state seen in the authored workbook.
$ws.Range("F25").Select() | Out-Null
